$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the four input cells that drive the whole calculator ---
# B2 = HbA1c (%), B3 = T1D duration (yrs), B4 = Incidence rate (%)
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 7.3
$ws.Range("B4").Value = 1.5

# B5 = Starting / Reference date. The sheet is protected and this cell is
# currently locked, so unlock it (matches the new data validation /
# unlocked-protection style that ships with this edit) before writing.
$ws.Range("B5").Locked = $false
$ws.Range("B5").Value = 44445

# --- View state: the author scrolled/zoomed and selected B2 before saving ---
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 240
